{"js": "// Insert a new citation paragraph before the first (and only) paragraph\n// in the document body: \"Nao ha certezas, apenas oportunidades. (V de Vingaca)\"\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = body.paragraphs.items[0];\nfirstParagraph.insertParagraph(\n  \"Nao ha certezas, apenas oportunidades. (V de Vingaca)\",\n  \"Before\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new citation paragraph before the first (and only) paragraph\n# in the document body: \"Nao ha certezas, apenas oportunidades. (V de Vingaca)\"\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs(1)\n$firstParagraph.Range.InsertParagraphBefore()\n\n$newParagraph = $d.Paragraphs(1)\n$newParagraph.Range.Text = \"Nao ha certezas, apenas oportunidades. (V de Vingaca)\"\n"}
